$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new volunteer-hour entries as rows 30 and 31 (pushing the previous
# "Total Project Hours:" summary row's formula range to include them).
$ws.Range("A30").Value = "11:31AM 2-25-2018"
$ws.Range("B30").Value = "1:17PM 2-25-2018"
$ws.Range("C30").Value = 106

$ws.Range("A31").Value = "4:43PM 2-25-2018"
$ws.Range("B31").Value = "9:11PM 2-25-2018"
$ws.Range("C31").Value = 268

# Update the selected cell to reflect the new last-entered data cell.
$ws.Range("C31").Select()
